$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Stage-2 bullet: " But, we have to find the output employee ..."
#    becomes " But, we find the output as employee ..."
#    Done as two narrow replacements so we touch only the words that
#    actually change and leave the rest of the sentence alone.
# ------------------------------------------------------------------
$d.Content.Find.Execute("e have to", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e", 2) | Out-Null

$d.Content.Find.Execute("the output employee", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the output as employee", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Semi Supervised Learning" -> "Supervised Learning"
#    (appears twice: once in the Stage-2 paragraph, once in the
#    Machine Learning / ... / Classification diagram line)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Semi Supervised Learning", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Supervised Learning", 2) | Out-Null

# ------------------------------------------------------------------
# 3) The small arrow autoshape that sits over "Supervised Learning"
#    shifts left (285.45pt -> 249.6pt) once the label text shortens.
#    It is a legacy VML drawing (<w:pict>/<v:shape>), which this
#    document's object model does not surface through
#    Document.Shapes/InlineShapes, so it cannot be reached from here.
#    Best-effort guarded attempt in case a shape handle ever becomes
#    available; silently skipped otherwise.
# ------------------------------------------------------------------
try {
    for ($i = 1; $i -le $d.Shapes.Count; $i++) {
        $shp = $d.Shapes.Item($i)
        if ($shp -ne $null -and $shp.Left -gt 280 -and $shp.Left -lt 290) {
            $shp.Left = 249.6
        }
    }
} catch {
    # Legacy VML shape not addressable via the Shapes collection here; ignore.
}
